$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hlookup")
$ws.Activate()

# Row 11 ("dmart"): HLOOKUP with exact-match 4th arg omitted, row index 2
$ws.Range("E11").Formula = '=HLOOKUP(E10,$D$20:$J$22,2,)'
$ws.Range("F11:J11").Formula = '=HLOOKUP(F10,$D$20:$J$22,2,)'

# Row 14 ("swiggy instamart"): HLOOKUP with explicit FALSE, row index 3
$ws.Range("E14").Formula = '=HLOOKUP(E10,$D$20:$J$22,3,FALSE)'
$ws.Range("F14").Formula = '=HLOOKUP(F10,$D$20:$J$22,3,FALSE)'
$ws.Range("G14:J14").Formula = '=HLOOKUP(G10,$D$20:$J$22,3,FALSE)'

$ws.Range("E27").Select()
